$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.515.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.624.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "203.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +10.51%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "568.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.619.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.623"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.678"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +16.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.01%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +11.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.193.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.606.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.04%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.343.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "405.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +16.78%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +10.32%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +16.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "674.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.10%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.423"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.96%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.280.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.24%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +10.37%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +29.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0421"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +12.13%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.11%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.81%  "
$ws.Range("E51").Style = "Normal"
